$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "The Other Side" paragraph right after the "SmartCard"
#    paragraph that precedes "Resources".
# ---------------------------------------------------------------------------
$rngSmartCard = $d.Content
$null = $rngSmartCard.Find.Execute("SmartCard", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pSmartCard = $rngSmartCard.Paragraphs(1)
$null = $pSmartCard.Range.InsertParagraphAfter()

$rngAfterSmartCard = $d.Content
$null = $rngAfterSmartCard.Find.Execute("SmartCard", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pOtherSideTarget = $rngAfterSmartCard.Paragraphs(1).Next()

$xmlOtherSide = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">The Other Side</w:t>
  </w:r>
</w:p>
"@
$null = $pOtherSideTarget.Range.InsertXML($xmlOtherSide)

# ---------------------------------------------------------------------------
# 2) After the "Exchanges" paragraph, insert two new paragraphs:
#      "Exchanges Listing Guide"  (ind left=360)
#      "Graphics"                  (ind left=360 firstLine=360, duplicate)
#    Then rewrite the original "Graphics" paragraph's run so it reads
#    "Marketing Materials" in the Times New Roman style (tab + text).
# ---------------------------------------------------------------------------
$rngExchanges = $d.Content
$null = $rngExchanges.Find.Execute("Exchanges", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pExchanges = $rngExchanges.Paragraphs(1)
$null = $pExchanges.Range.InsertParagraphAfter()

$rngAfterExchanges = $d.Content
$null = $rngAfterExchanges.Find.Execute("Exchanges", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pListingGuideTarget = $rngAfterExchanges.Paragraphs(1).Next()

$xmlListingGuide = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:tab/>
    <w:t xml:space="preserve">Exchanges Listing Guide</w:t>
  </w:r>
</w:p>
"@
$null = $pListingGuideTarget.Range.InsertXML($xmlListingGuide)

$null = $pListingGuideTarget.Range.InsertParagraphAfter()

$rngListingGuide = $d.Content
$null = $rngListingGuide.Find.Execute("Exchanges Listing Guide", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pGraphicsDupTarget = $rngListingGuide.Paragraphs(1).Next()

$xmlGraphicsDup = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360" w:firstLine="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Graphics</w:t>
  </w:r>
</w:p>
"@
$null = $pGraphicsDupTarget.Range.InsertXML($xmlGraphicsDup)

# Now find the ORIGINAL "Graphics" paragraph (the one that still has the
# firstLine-indent pPr with no preceding tab) and rewrite its run.
$rngOldGraphics = $d.Content
$null = $rngOldGraphics.Find.Execute("Graphics", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$null = $rngOldGraphics.Find.Execute("Graphics", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pOldGraphics = $rngOldGraphics.Paragraphs(1)

$xmlMarketingMaterials = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360" w:firstLine="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">Marketing Materials</w:t>
  </w:r>
</w:p>
"@
$null = $pOldGraphics.Range.InsertXML($xmlMarketingMaterials)

# ---------------------------------------------------------------------------
# 3) Insert a new "Guide" paragraph right after "Electrum Wallet".
# ---------------------------------------------------------------------------
$rngElectrum = $d.Content
$null = $rngElectrum.Find.Execute("Electrum Wallet", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pElectrum = $rngElectrum.Paragraphs(1)
$null = $pElectrum.Range.InsertParagraphAfter()

$rngAfterElectrum = $d.Content
$null = $rngAfterElectrum.Find.Execute("Electrum Wallet", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pGuideTarget = $rngAfterElectrum.Paragraphs(1).Next()

$xmlGuide = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360" w:firstLine="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:tab/>
    <w:t>Guide</w:t>
  </w:r>
</w:p>
"@
$null = $pGuideTarget.Range.InsertXML($xmlGuide)

Write-Host "done"
